# Apply the commit: insert 4 new weekly rows (Ajo / Chino) right before the
# existing row 879 block, pushing the remaining historical rows down by 4
# (old row N -> new row N+4, for N = 879..954).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows at row 879 (each Insert() pushes the row at that
# index - and everything below it - down by one, Excel-style).
$ws.Rows.Item(879).Insert()
$ws.Rows.Item(879).Insert()
$ws.Rows.Item(879).Insert()
$ws.Rows.Item(879).Insert()

# Common columns shared by every data row in this sheet.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112003
$categoria = "Ajo"
$clasificacion = "Hortaliza"

function Set-AjoRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $categoriaId
    $ws.Cells.Item($Row, 7).Value = $categoria
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasificacion
}

Set-AjoRow 879 44931 "Chino" "1a (cosecha)" 500 12000 12000 12000 "$/caja 10 kilos" "Región Metropolitana" 1200 10
Set-AjoRow 880 44931 "Chino" "2a (cosecha)" 300 11000 11000 11000 "$/caja 10 kilos" "Región Metropolitana" 1100 10
Set-AjoRow 881 44931 "Chino" "3a (cosecha)" 100 10000 10000 10000 "$/caja 10 kilos" "Región Metropolitana" 1000 10
Set-AjoRow 882 44931 "Chino" "Primera"      1100 12000 12500 12273 "$/caja 10 kilos" "China" 1227 10
